$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1133093333333333
$ws.Range("H2").Value = 0.339928
$ws.Range("I2").Value = 0.02456654176752224
$ws.Range("J2").Value = 0.02456654176752224
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.070065
$ws.Range("N2").Value = 3.210195
$ws.Range("O2").Value = 0.07699366399823478
$ws.Range("P2").Value = 0.07699366399823478
$ws.Range("Q2").Value = 0.1212483517733333
$ws.Range("R2").Value = 1.09123516596
$ws.Range("S2").Value = 0.001891468062447208
$ws.Range("T2").Value = 0.001891468062447208
$ws.Range("G3").Value = 0.1133093333333333
$ws.Range("H3").Value = 0.339928
$ws.Range("I3").Value = 0.02456654176752224
$ws.Range("J3").Value = 0.02456654176752224
$ws.Range("O3").Value = 0.08667273864337491
$ws.Range("P3").Value = 0.08667273864337491
$ws.Range("Q3").Value = 0.1364907988328889
$ws.Range("R3").Value = 1.228417189496
$ws.Range("S3").Value = 0.002129249453988008
$ws.Range("T3").Value = 0.002129249453988008
$ws.Range("G4").Value = 0.1133093333333333
$ws.Range("H4").Value = 0.339928
$ws.Range("I4").Value = 0.02456654176752224
$ws.Range("J4").Value = 0.02456654176752224
$ws.Range("M4").Value = 11.623441
$ws.Range("N4").Value = 34.870323
$ws.Range("O4").Value = 0.8363335973583904
$ws.Range("P4").Value = 0.8363335973583904
$ws.Range("Q4").Value = 1.317044350749333
$ws.Range("R4").Value = 11.853399156744
$ws.Range("S4").Value = 0.02054582425108702
$ws.Range("T4").Value = 0.02054582425108702
$ws.Range("I5").Value = 0.8380577451911468
$ws.Range("J5").Value = 0.8380577451911468
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.070065
$ws.Range("N5").Value = 3.210195
$ws.Range("O5").Value = 0.07699366399823478
$ws.Range("P5").Value = 0.07699366399823478
$ws.Range("Q5").Value = 4.136240308338333
$ws.Range("R5").Value = 37.226162775045
$ws.Range("S5").Value = 0.06452513644436542
$ws.Range("T5").Value = 0.06452513644436542
$ws.Range("I6").Value = 0.8380577451911468
$ws.Range("J6").Value = 0.8380577451911468
$ws.Range("O6").Value = 0.08667273864337491
$ws.Range("P6").Value = 0.08667273864337491
$ws.Range("S6").Value = 0.07263675991700835
$ws.Range("T6").Value = 0.07263675991700835
$ws.Range("I7").Value = 0.8380577451911468
$ws.Range("J7").Value = 0.8380577451911468
$ws.Range("M7").Value = 11.623441
$ws.Range("N7").Value = 34.870323
$ws.Range("O7").Value = 0.8363335973583904
$ws.Range("P7").Value = 0.8363335973583904
$ws.Range("Q7").Value = 44.92936895029033
$ws.Range("R7").Value = 404.364320552613
$ws.Range("S7").Value = 0.7008958488297732
$ws.Range("T7").Value = 0.7008958488297732
$ws.Range("G8").Value = 0.6336240000000001
$ws.Range("H8").Value = 1.900872
$ws.Range("I8").Value = 0.1373757130413309
$ws.Range("J8").Value = 0.1373757130413309
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.070065
$ws.Range("N8").Value = 3.210195
$ws.Range("O8").Value = 0.07699366399823478
$ws.Range("P8").Value = 0.07699366399823478
$ws.Range("Q8").Value = 0.6780188655600001
$ws.Range("R8").Value = 6.10216979004
$ws.Range("S8").Value = 0.01057705949142215
$ws.Range("T8").Value = 0.01057705949142215
$ws.Range("G9").Value = 0.6336240000000001
$ws.Range("H9").Value = 1.900872
$ws.Range("I9").Value = 0.1373757130413309
$ws.Range("J9").Value = 0.1373757130413309
$ws.Range("O9").Value = 0.08667273864337491
$ws.Range("P9").Value = 0.08667273864337491
$ws.Range("Q9").Value = 0.7632543884560001
$ws.Range("R9").Value = 6.869289496104001
$ws.Range("S9").Value = 0.01190672927237854
$ws.Range("T9").Value = 0.01190672927237854
$ws.Range("G10").Value = 0.6336240000000001
$ws.Range("H10").Value = 1.900872
$ws.Range("I10").Value = 0.1373757130413309
$ws.Range("J10").Value = 0.1373757130413309
$ws.Range("M10").Value = 11.623441
$ws.Range("N10").Value = 34.870323
$ws.Range("O10").Value = 0.8363335973583904
$ws.Range("P10").Value = 0.8363335973583904
$ws.Range("Q10").Value = 7.364891180184001
$ws.Range("R10").Value = 66.284020621656
$ws.Range("S10").Value = 0.1148919242775302
$ws.Range("T10").Value = 0.1148919242775302
